$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that were removed from the export (bottom-up so row
# numbers above the deletion point stay valid):
#   row 9  -> 004643153 / CARLA    / 3800
#   row 8  -> 005890232 / TAYLA    / 4227
#   row 7  -> 001651617 / MIRELLA  / 4737.39
#   row 6  -> 003553997 / MIRELLA  / 6177.71
#   row 4  -> 004332747 / LOHRAN   / 6385.18
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()

# After the deletions above, the row that used to hold
# 004995535 / ASIEL / 2540.34 is now row 5. Replace its contents with the
# new record (keep the account number as text so the leading zero survives).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "004487140"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "VALMIR"
$ws.Range("C5").Value = 3270
